$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the affected data columns (Price, Volume(1h),
# Hora) so that numeric-looking strings (prices), percentages, and the plain
# digit in column G are stored as literal text, matching the original
# inline-string cell types instead of being auto-coerced into Excel
# Number/Percentage types by COM's input parser. Column F (Data/date) is
# intentionally left untouched since it is not modified by this edit.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "307.84"
$ws.Range("E2").Value = "0.27%"
$ws.Range("G2").Value = "4"

# Row 3
$ws.Range("D3").Value = "41.06"
$ws.Range("E3").Value = "3.05%"
$ws.Range("G3").Value = "4"

# Row 4
$ws.Range("D4").Value = "5.138"
$ws.Range("E4").Value = "2.40%"
$ws.Range("G4").Value = "4"

# Row 5
$ws.Range("D5").Value = "0.07615"
$ws.Range("E5").Value = "-0.72%"
$ws.Range("G5").Value = "4"

# Row 6
$ws.Range("D6").Value = "1.624"
$ws.Range("E6").Value = "0.34%"
$ws.Range("G6").Value = "4"

# Row 7
$ws.Range("D7").Value = "2.496"
$ws.Range("E7").Value = "-0.35%"
$ws.Range("G7").Value = "4"

# Row 8
$ws.Range("D8").Value = "0.9012"
$ws.Range("E8").Value = "1.84%"
$ws.Range("G8").Value = "4"

# Row 9
$ws.Range("D9").Value = "0.1119"
$ws.Range("E9").Value = "11.22%"
$ws.Range("G9").Value = "4"

# Row 10
$ws.Range("D10").Value = "0.1773"
$ws.Range("E10").Value = "2.51%"
$ws.Range("G10").Value = "4"

# Row 11
$ws.Range("D11").Value = "0.09246"
$ws.Range("E11").Value = "3.53%"
$ws.Range("G11").Value = "4"

# Row 12
$ws.Range("D12").Value = "0.04182"
$ws.Range("E12").Value = "-4.71%"
$ws.Range("G12").Value = "4"

# Row 13
$ws.Range("D13").Value = "0.1049"
$ws.Range("E13").Value = "-0.54%"
$ws.Range("G13").Value = "4"

# Row 14
$ws.Range("D14").Value = "0.001249"
$ws.Range("E14").Value = "-2.13%"
$ws.Range("G14").Value = "4"

# Row 15
$ws.Range("D15").Value = "0.005897"
$ws.Range("E15").Value = "1.61%"
$ws.Range("G15").Value = "4"

# Row 16
$ws.Range("D16").Value = "3.356"
$ws.Range("E16").Value = "0.12%"
$ws.Range("G16").Value = "4"

# Row 17
$ws.Range("D17").Value = "4.236"
$ws.Range("E17").Value = "-0.21%"
$ws.Range("G17").Value = "4"

# Row 18
$ws.Range("E18").Value = "-1.94%"
$ws.Range("G18").Value = "4"

# Row 19
$ws.Range("D19").Value = "6.560"
$ws.Range("E19").Value = "-6.19%"
$ws.Range("G19").Value = "4"

# Row 20
$ws.Range("D20").Value = "0.1364"
$ws.Range("E20").Value = "1.57%"
$ws.Range("G20").Value = "4"

# Row 21
$ws.Range("D21").Value = "0.2681"
$ws.Range("E21").Value = "-18.96%"
$ws.Range("G21").Value = "4"

# Row 22
$ws.Range("D22").Value = "0.04145"
$ws.Range("E22").Value = "-1.63%"
$ws.Range("G22").Value = "4"

# Row 23
$ws.Range("D23").Value = "0.001228"
$ws.Range("E23").Value = "2.26%"
$ws.Range("G23").Value = "4"

# Row 24
$ws.Range("D24").Value = "0.004080"
$ws.Range("E24").Value = "0.51%"
$ws.Range("G24").Value = "4"

# Row 25
$ws.Range("E25").Value = "6.28%"
$ws.Range("G25").Value = "4"

# Row 26
$ws.Range("G26").Value = "4"

# Row 27
$ws.Range("G27").Value = "4"

# Row 28
$ws.Range("G28").Value = "4"

# Row 29
$ws.Range("G29").Value = "4"

# Row 30
$ws.Range("G30").Value = "4"

# Row 31
$ws.Range("G31").Value = "4"

# Row 32
$ws.Range("G32").Value = "4"

# Row 33
$ws.Range("G33").Value = "4"

# Row 34
$ws.Range("G34").Value = "4"

# Row 35
$ws.Range("G35").Value = "4"

# Row 36
$ws.Range("G36").Value = "4"

# Row 37
$ws.Range("G37").Value = "4"

# Row 38
$ws.Range("D38").Value = "0.02407"
$ws.Range("E38").Value = "2.63%"
$ws.Range("G38").Value = "4"

# Row 39
$ws.Range("D39").Value = "0.05190"
$ws.Range("E39").Value = "0.35%"
$ws.Range("G39").Value = "4"

# Row 40
$ws.Range("D40").Value = "0.007771"
$ws.Range("E40").Value = "-2.46%"
$ws.Range("G40").Value = "4"

# Row 41
$ws.Range("D41").Value = "0.1301"
$ws.Range("E41").Value = "-1.54%"
$ws.Range("G41").Value = "4"

# Row 42
$ws.Range("D42").Value = "0.006959"
$ws.Range("E42").Value = "5.73%"
$ws.Range("G42").Value = "4"

# Row 43
$ws.Range("D43").Value = "0.001970"
$ws.Range("E43").Value = "-1.38%"
$ws.Range("G43").Value = "4"

# Row 44
$ws.Range("D44").Value = "0.007573"
$ws.Range("E44").Value = "-1.08%"
$ws.Range("G44").Value = "4"

# Row 45
$ws.Range("D45").Value = "0.3052"
$ws.Range("E45").Value = "0.22%"
$ws.Range("G45").Value = "4"

# Row 46
$ws.Range("D46").Value = "0.00006733"
$ws.Range("E46").Value = "2.25%"
$ws.Range("G46").Value = "4"

# Row 47
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "-0.18%"
$ws.Range("G47").Value = "4"

# Row 48
$ws.Range("D48").Value = "0.03165"
$ws.Range("E48").Value = "835.12%"
$ws.Range("G48").Value = "4"

# Row 49
$ws.Range("E49").Value = "-15.98%"
$ws.Range("G49").Value = "4"

# Row 50
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "-0.18%"
$ws.Range("G50").Value = "4"

# Row 51
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "-0.18%"
$ws.Range("G51").Value = "4"
